# Refatorando o consolidador para modelo ETL
# Updates the absenteeism data rows (2-11) with new source values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 65446
$ws.Cells.Item(2, 2).Value = "Breno Freitas"
$ws.Cells.Item(2, 3).Value = "P&D"
$ws.Cells.Item(2, 4).Value = "Consulta médica"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 45086
$ws.Cells.Item(2, 7).Value = 7302.3

# Row 3
$ws.Cells.Item(3, 1).Value = 76922
$ws.Cells.Item(3, 2).Value = "João Vitor Porto"
$ws.Cells.Item(3, 3).Value = "Recursos Humanos"
$ws.Cells.Item(3, 4).Value = "Consulta médica"
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = 45082
$ws.Cells.Item(3, 7).Value = 12311.28

# Row 4
$ws.Cells.Item(4, 1).Value = 77723
$ws.Cells.Item(4, 2).Value = "Julia Pereira"
$ws.Cells.Item(4, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(4, 4).Value = "Outros"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 45100
$ws.Cells.Item(4, 7).Value = 9653.98

# Row 5
$ws.Cells.Item(5, 1).Value = 12199
$ws.Cells.Item(5, 2).Value = "Luana Correia"
$ws.Cells.Item(5, 3).Value = "Engenharia"
$ws.Cells.Item(5, 4).Value = "Doença"
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 45095
$ws.Cells.Item(5, 7).Value = 9501.639999999999

# Row 6
$ws.Cells.Item(6, 1).Value = 5228
$ws.Cells.Item(6, 2).Value = "Ana Sophia Moraes"
$ws.Cells.Item(6, 3).Value = "Vendas"
$ws.Cells.Item(6, 4).Value = "Doença"
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 45100
$ws.Cells.Item(6, 7).Value = 5171.68

# Row 7
$ws.Cells.Item(7, 1).Value = 86422
$ws.Cells.Item(7, 2).Value = "João Lucas Fernandes"
$ws.Cells.Item(7, 3).Value = "P&D"
$ws.Cells.Item(7, 4).Value = "Outros"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 45078
$ws.Cells.Item(7, 7).Value = 9726.26

# Row 8
$ws.Cells.Item(8, 1).Value = 70845
$ws.Cells.Item(8, 2).Value = "Helena da Mata"
$ws.Cells.Item(8, 3).Value = "Marketing"
$ws.Cells.Item(8, 4).Value = "Outros"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 45089
$ws.Cells.Item(8, 7).Value = 9203.93

# Row 9
$ws.Cells.Item(9, 1).Value = 93430
$ws.Cells.Item(9, 2).Value = "Carlos Eduardo Moraes"
$ws.Cells.Item(9, 3).Value = "Engenharia"
$ws.Cells.Item(9, 4).Value = "Problemas pessoais"
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 45094
$ws.Cells.Item(9, 7).Value = 8471.74

# Row 10
$ws.Cells.Item(10, 1).Value = 97545
$ws.Cells.Item(10, 2).Value = "Laura Rezende"
$ws.Cells.Item(10, 3).Value = "Jurídico"
$ws.Cells.Item(10, 4).Value = "Problemas pessoais"
$ws.Cells.Item(10, 5).Value = 6
$ws.Cells.Item(10, 6).Value = 45083
$ws.Cells.Item(10, 7).Value = 9568.4

# Row 11
$ws.Cells.Item(11, 1).Value = 35560
$ws.Cells.Item(11, 2).Value = "Clarice Moreira"
$ws.Cells.Item(11, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(11, 4).Value = "Problemas pessoais"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 45079
$ws.Cells.Item(11, 7).Value = 10253.07
